$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold text-formatted numbers/percentages
# in the source data; force text format so Excel does not auto-convert
# these assignments into numeric cell values.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '67.508.79'
$ws.Range('E2').Value = '  -2.48%  '
$ws.Range('D3').Value = '3.261.12'
$ws.Range('E3').Value = '  -5.28%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '590.72'
$ws.Range('E5').Value = '  -3.05%  '
$ws.Range('D6').Value = '149.53'
$ws.Range('E6').Value = '  -10.73%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.254.07'
$ws.Range('E8').Value = '  -5.33%  '
$ws.Range('D9').Value = '0.543'
$ws.Range('E9').Value = '  -8.77%  '
$ws.Range('E10').Value = '  -11.54%  '
$ws.Range('D11').Value = '6.72'
$ws.Range('E11').Value = '  -4.69%  '
$ws.Range('D12').Value = '0.504'
$ws.Range('E12').Value = '  -10.48%  '
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').Value = '38.37'
$ws.Range('E13').Value = '  -13.21%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').Value = '0.0000246'
$ws.Range('E14').Value = '  -8.73%  '
$ws.Range('D15').Value = '3.781.13'
$ws.Range('E15').Value = '  -5.41%  '
$ws.Range('D16').Value = '67.508.82'
$ws.Range('E16').Value = '  -2.60%  '
$ws.Range('D17').Value = '3.259.58'
$ws.Range('E17').Value = '  -5.31%  '
$ws.Range('E18').Value = '  -5.48%  '
$ws.Range('D19').Value = '529.12'
$ws.Range('E19').Value = '  -8.91%  '
$ws.Range('E20').Value = '  -12.84%  '
$ws.Range('D21').Value = '14.95'
$ws.Range('E21').Value = '  -12.68%  '
$ws.Range('D22').Value = '0.754'
$ws.Range('E22').Value = '  -10.82%  '
$ws.Range('D23').Value = '7.85'
$ws.Range('E23').Value = '  -11.87%  '
$ws.Range('D24').Value = '85.57'
$ws.Range('E24').Value = '  -11.03%  '
$ws.Range('D25').Value = '13.48'
$ws.Range('E25').Value = '  -11.01%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').Value = '3.20'
$ws.Range('E27').Value = '  -11.24%  '
$ws.Range('E28').Value = '  -11.74%  '
$ws.Range('D29').Value = '8.00'
$ws.Range('E29').Value = '  -7.30%  '
$ws.Range('D30').Value = '28.97'
$ws.Range('E30').Value = '  -11.54%  '
$ws.Range('E31').Value = '  -3.17%  '
$ws.Range('E32').Value = '  -4.69%  '
$ws.Range('D33').Value = '6.60'
$ws.Range('E33').Value = '  -15.55%  '
$ws.Range('E34').Value = '  -13.31%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Value = '512.96'
$ws.Range('E36').Value = '  -11.84%  '
$ws.Range('D37').Value = '0.0442'
$ws.Range('E37').Value = '  -6.56%  '
$ws.Range('D38').Value = '52.97'
$ws.Range('E38').Value = '  -5.50%  '
$ws.Range('D39').Value = '0.0852'
$ws.Range('E39').Value = '  -10.64%  '
$ws.Range('D40').Value = '8.92'
$ws.Range('E40').Value = '  -15.07%  '
$ws.Range('E41').Value = '  -10.94%  '
$ws.Range('D42').Value = '2.78'
$ws.Range('E42').Value = '  -11.46%  '
$ws.Range('D43').Value = '2.931.22'
$ws.Range('E43').Value = '  -9.47%  '
$ws.Range('D44').Value = '0.266'
$ws.Range('E44').Value = '  -9.92%  '
$ws.Range('D45').Value = '0.0₃0588'
$ws.Range('E45').Value = '  -14.31%  '
$ws.Range('D46').Value = '2.19'
$ws.Range('E46').Value = '  -8.70%  '
$ws.Range('D47').Value = '26.55'
$ws.Range('E47').Value = '  -14.44%  '
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('D49').Value = '2.31'
$ws.Range('E49').Value = '  -16.78%  '
$ws.Range('E50').Value = '  -10.26%  '
$ws.Range('D51').Value = '123.73'
$ws.Range('E51').Value = '  -7.62%  '
